# "push dynamic sanksi with multiple table page"
#
# The sample sanksi row (row 2) advances to the next NPP/Kode_PKS value:
#   A2: "01733722" -> "01733725"
#
# A2 uses a text-quoted style (quotePrefix / numFmtId 49) so the leading
# zeros in the numeric-looking code are preserved as text. A leading
# apostrophe keeps Excel treating the entry as quoted text (instead of
# re-coercing the cell to a generic text style), so the cell's existing
# formatting stays intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'01733725"
